$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebased bus-voltage results (vm_pu) for the 380 kV slack-bus case (Case_5_234).
# Each assignment below mirrors one cell-value edit from the recorded diff.
# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.070955003893758
$ws.Range("D2").Value = 1.076060471359031
$ws.Range("E2").Value = 1.065108339823106
$ws.Range("F2").Value = 1.08496578916929
$ws.Range("I2").Value = 1.062660197921192
$ws.Range("J2").Value = 1.075882042693381
$ws.Range("K2").Value = 1.078745398953193
$ws.Range("L2").Value = 1.067822468811518
$ws.Range("M2").Value = 1.087627441953888
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.072325022927828
$ws.Range("D3").Value = 1.077198851535017
$ws.Range("E3").Value = 1.066316203070301
$ws.Range("F3").Value = 1.086236525146257
$ws.Range("I3").Value = 1.063201114102707
$ws.Range("J3").Value = 1.076907821865162
$ws.Range("K3").Value = 1.079700426395129
$ws.Range("L3").Value = 1.068844685466318
$ws.Range("M3").Value = 1.088716194643674
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.073210065510805
$ws.Range("D4").Value = 1.077934222002555
$ws.Range("E4").Value = 1.067095971961824
$ws.Range("F4").Value = 1.087057808314186
$ws.Range("I4").Value = 1.063549131949797
$ws.Range("J4").Value = 1.077569620331845
$ws.Range("K4").Value = 1.08031656859472
$ws.Range("L4").Value = 1.069503760708313
$ws.Range("M4").Value = 1.089419159242752
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.073581795873477
$ws.Range("D5").Value = 1.078243079671593
$ws.Range("E5").Value = 1.067423361206658
$ws.Range("F5").Value = 1.087402848627433
$ws.Range("I5").Value = 1.063694964593371
$ws.Range("J5").Value = 1.077847378251785
$ws.Range("K5").Value = 1.080575162022377
$ws.Range("L5").Value = 1.069780274036464
$ws.Range("M5").Value = 1.089714323291805
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.073644191180414
$ws.Range("D6").Value = 1.078294921227862
$ws.Range("E6").Value = 1.067478306494346
$ws.Range("F6").Value = 1.087460769228626
$ws.Range("I6").Value = 1.063719422812496
$ws.Range("J6").Value = 1.077893988058017
$ws.Range("K6").Value = 1.080618555736189
$ws.Range("L6").Value = 1.069826669030599
$ws.Range("M6").Value = 1.089763861515428
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.073215033924772
$ws.Range("D7").Value = 1.077938350116273
$ws.Range("E7").Value = 1.067100348218802
$ws.Range("F7").Value = 1.087062419646487
$ws.Range("I7").Value = 1.063551082432213
$ws.Range("J7").Value = 1.07757333355714
$ws.Range("K7").Value = 1.080320025628202
$ws.Range("L7").Value = 1.069507457694836
$ws.Range("M7").Value = 1.089423104657352
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.071418312319541
$ws.Range("D8").Value = 1.076445451072665
$ws.Range("E8").Value = 1.065516918566539
$ws.Range("F8").Value = 1.085395443549037
$ws.Range("I8").Value = 1.062843416922727
$ws.Range("J8").Value = 1.076229115886992
$ws.Range("K8").Value = 1.07906853512987
$ws.Range("L8").Value = 1.068168424792613
$ws.Range("M8").Value = 1.087995710663569
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.068240877960937
$ws.Range("D9").Value = 1.073805103750211
$ws.Range("E9").Value = 1.062712699032934
$ws.Range("F9").Value = 1.082450397870894
$ws.Range("I9").Value = 1.061581052946452
$ws.Range("J9").Value = 1.073845301376758
$ws.Range("K9").Value = 1.076849095439672
$ws.Range("L9").Value = 1.065790530059979
$ws.Range("M9").Value = 1.085468539376887
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.066114562561246
$ws.Range("D10").Value = 1.072038104032136
$ws.Range("E10").Value = 1.060833469181777
$ws.Range("F10").Value = 1.080481618091814
$ws.Range("I10").Value = 1.060728979362916
$ws.Range("J10").Value = 1.072245633443683
$ws.Range("K10").Value = 1.075359698263004
$ws.Range("L10").Value = 1.064192625556583
$ws.Range("M10").Value = 1.08377547083169
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.065191857140124
$ws.Range("D11").Value = 1.0712713080548
$ws.Range("E11").Value = 1.060017355716842
$ws.Range("F11").Value = 1.079627764360899
$ws.Range("I11").Value = 1.060357494890997
$ws.Range("J11").Value = 1.071550415009401
$ws.Range("K11").Value = 1.074712398066148
$ws.Range("L11").Value = 1.063497648373692
$ws.Range("M11").Value = 1.083040326337599
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.064848815960367
$ws.Range("D12").Value = 1.070986229187464
$ws.Range("E12").Value = 1.059713848654372
$ws.Range("F12").Value = 1.079310395122609
$ws.Range("I12").Value = 1.060219125481587
$ws.Range("J12").Value = 1.071291790966703
$ws.Range("K12").Value = 1.074471599150723
$ws.Range("L12").Value = 1.063239035076239
$ws.Range("M12").Value = 1.082766949878216
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.064922413481803
$ws.Range("D13").Value = 1.071047391273715
$ws.Range("E13").Value = 1.059778968623904
$ws.Range("F13").Value = 1.079378481494884
$ws.Range("I13").Value = 1.060248823598408
$ws.Range("J13").Value = 1.071347284393503
$ws.Range("K13").Value = 1.074523267834078
$ws.Range("L13").Value = 1.063294529773327
$ws.Range("M13").Value = 1.082825604219769
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.065163507560496
$ws.Range("D14").Value = 1.071247748614
$ws.Range("E14").Value = 1.059992275235783
$ws.Range("F14").Value = 1.079601534836913
$ws.Range("I14").Value = 1.06034606507555
$ws.Range("J14").Value = 1.071529045032601
$ws.Range("K14").Value = 1.074692500982908
$ws.Range("L14").Value = 1.063476280896219
$ws.Range("M14").Value = 1.083017735331014
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.065312012639563
$ws.Range("D15").Value = 1.071371161210006
$ws.Range("E15").Value = 1.060123651763262
$ws.Range("F15").Value = 1.079738937369989
$ws.Range("I15").Value = 1.060405927838091
$ws.Range("J15").Value = 1.071640982051738
$ws.Range("K15").Value = 1.074796722885693
$ws.Range("L15").Value = 1.063588201589801
$ws.Range("M15").Value = 1.083136072269439
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.06617575678823
$ws.Range("D16").Value = 1.072088958048185
$ws.Range("E16").Value = 1.060887580946483
$ws.Range("F16").Value = 1.080538256426921
$ws.Range("I16").Value = 1.060753579959092
$ws.Range("J16").Value = 1.072291718563555
$ws.Range("K16").Value = 1.075402606844056
$ws.Range("L16").Value = 1.06423868362709
$ws.Range("M16").Value = 1.083824216624742
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.066717020926358
$ws.Range("D17").Value = 1.072538761295075
$ws.Range("E17").Value = 1.061366127927156
$ws.Range("F17").Value = 1.081039280234528
$ws.Range("I17").Value = 1.060970973021193
$ws.Range("J17").Value = 1.072699221288854
$ws.Range("K17").Value = 1.07578202074205
$ws.Range("L17").Value = 1.064645886608836
$ws.Range("M17").Value = 1.084255322880697
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.067032538972874
$ws.Range("D18").Value = 1.072800962816098
$ws.Range("E18").Value = 1.061645025683791
$ws.Range("F18").Value = 1.081331388203751
$ws.Range("I18").Value = 1.061097530736167
$ws.Range("J18").Value = 1.072936664916768
$ws.Range("K18").Value = 1.076003097021784
$ws.Range("L18").Value = 1.064883105031117
$ws.Range("M18").Value = 1.084506583918807
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.067140090097917
$ws.Range("D19").Value = 1.07289033959199
$ws.Range("E19").Value = 1.06174008369116
$ws.Range("F19").Value = 1.081430967469713
$ws.Range("I19").Value = 1.061140642361984
$ws.Range("J19").Value = 1.07301758553155
$ws.Range("K19").Value = 1.076078439556293
$ws.Range("L19").Value = 1.064963940324808
$ws.Range("M19").Value = 1.084592224365453
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.066658968338159
$ws.Range("D20").Value = 1.072490518366787
$ws.Range("E20").Value = 1.061314808265241
$ws.Range("F20").Value = 1.080985538687892
$ws.Range("I20").Value = 1.060947674063613
$ws.Range("J20").Value = 1.072655525570932
$ws.Range("K20").Value = 1.075741336989926
$ws.Range("L20").Value = 1.064602228263493
$ws.Range("M20").Value = 1.084209089567402
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.065092519905396
$ws.Range("D21").Value = 1.071188755499872
$ws.Range("E21").Value = 1.059929471915463
$ws.Range("F21").Value = 1.07953585702909
$ws.Range("I21").Value = 1.060317440498178
$ws.Range("J21").Value = 1.071475531837091
$ws.Range("K21").Value = 1.074642676094226
$ws.Range("L21").Value = 1.063422772683248
$ws.Range("M21").Value = 1.082961166137732
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.064105849991591
$ws.Range("D22").Value = 1.070368797687289
$ws.Range("E22").Value = 1.059056334123165
$ws.Range("F22").Value = 1.078623167075514
$ws.Range("I22").Value = 1.059918966976644
$ws.Range("J22").Value = 1.0707313697279
$ws.Range("K22").Value = 1.073949802812785
$ws.Range("L22").Value = 1.062678492441738
$ws.Range("M22").Value = 1.082174744746364
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.06462907370791
$ws.Range("D23").Value = 1.07080361564111
$ws.Range("E23").Value = 1.059519404271736
$ws.Range("F23").Value = 1.079107118631064
$ws.Range("I23").Value = 1.060130416925206
$ws.Range("J23").Value = 1.071126079601026
$ws.Range("K23").Value = 1.074317308939809
$ws.Range("L23").Value = 1.063073308340392
$ws.Range("M23").Value = 1.082591814258136
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.066685200397161
$ws.Range("D24").Value = 1.072512317766929
$ws.Range("E24").Value = 1.061337998123502
$ws.Range("F24").Value = 1.081009822583139
$ws.Range("I24").Value = 1.060958202613276
$ws.Range("J24").Value = 1.07267527054414
$ws.Range("K24").Value = 1.075759720934847
$ws.Range("L24").Value = 1.064621956505239
$ws.Range("M24").Value = 1.084229981016674
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.069063707975156
$ws.Range("D25").Value = 1.074488870017843
$ws.Range("E25").Value = 1.063439351369603
$ws.Range("F25").Value = 1.083212695830281
$ws.Range("I25").Value = 1.061909242346292
$ws.Range("J25").Value = 1.074463396906063
$ws.Range("K25").Value = 1.077424577174966
$ws.Range("L25").Value = 1.066407478462817
$ws.Range("M25").Value = 1.08612331436963

Write-Host "Updated vm_pu values for rows 2-25"
